{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces the 25 arithmetic-answer text cells inside the single table in\n// the document body with their new values, cell-by-cell, while preserving\n// each cell's existing run/paragraph formatting (font, size, alignment).\n//\n// Each entry below records the (row, col) of a content-bearing table cell\n// (0-based, matching Table.getCell), the text currently in that cell, and\n// the text it must become, taken straight from the authoritative diff.\nconst replacements = [\n  {\n    \"row\": 0,\n    \"col\": 0,\n    \"oldText\": \"90\u00f73=30, 0\",\n    \"newText\": \"74\u00f78=9, 2\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 1,\n    \"oldText\": \"75\u00f72=37, 1\",\n    \"newText\": \"79\u00f73=26, 1\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 2,\n    \"oldText\": \"15\u00f79=1, 6\",\n    \"newText\": \"31\u00f76=5, 1\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 3,\n    \"oldText\": \"76\u00f74=19, 0\",\n    \"newText\": \"40\u00f75=8, 0\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 4,\n    \"oldText\": \"10\u00f73=3, 1\",\n    \"newText\": \"98\u00f74=24, 2\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 0,\n    \"oldText\": \"74\u00f78=9, 2\",\n    \"newText\": \"48\u00f74=12, 0\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 1,\n    \"oldText\": \"86\u00f73=28, 2\",\n    \"newText\": \"40\u00f74=10, 0\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 2,\n    \"oldText\": \"35\u00f75=7, 0\",\n    \"newText\": \"66\u00f74=16, 2\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 3,\n    \"oldText\": \"10\u00f76=1, 4\",\n    \"newText\": \"59\u00f77=8, 3\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 4,\n    \"oldText\": \"23\u00f75=4, 3\",\n    \"newText\": \"95\u00f75=19, 0\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 0,\n    \"oldText\": \"27\u00f72=13, 1\",\n    \"newText\": \"91\u00f73=30, 1\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 1,\n    \"oldText\": \"86\u00f78=10, 6\",\n    \"newText\": \"98\u00f72=49, 0\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 2,\n    \"oldText\": \"45\u00f76=7, 3\",\n    \"newText\": \"15\u00f78=1, 7\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 3,\n    \"oldText\": \"66\u00f74=16, 2\",\n    \"newText\": \"36\u00f79=4, 0\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 4,\n    \"oldText\": \"57\u00f75=11, 2\",\n    \"newText\": \"56\u00f74=14, 0\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 0,\n    \"oldText\": \"17\u00f79=1, 8\",\n    \"newText\": \"31\u00f75=6, 1\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 1,\n    \"oldText\": \"11\u00f74=2, 3\",\n    \"newText\": \"88\u00f76=14, 4\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 2,\n    \"oldText\": \"90\u00f72=45, 0\",\n    \"newText\": \"15\u00f73=5, 0\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 3,\n    \"oldText\": \"73\u00f73=24, 1\",\n    \"newText\": \"31\u00f72=15, 1\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 4,\n    \"oldText\": \"48\u00f73=16, 0\",\n    \"newText\": \"29\u00f74=7, 1\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 0,\n    \"oldText\": \"81\u00f79=9, 0\",\n    \"newText\": \"76\u00f72=38, 0\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 1,\n    \"oldText\": \"55\u00f77=7, 6\",\n    \"newText\": \"61\u00f75=12, 1\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 2,\n    \"oldText\": \"90\u00f75=18, 0\",\n    \"newText\": \"62\u00f75=12, 2\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 3,\n    \"oldText\": \"84\u00f75=16, 4\",\n    \"newText\": \"98\u00f72=49, 0\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 4,\n    \"oldText\": \"80\u00f78=10, 0\",\n    \"newText\": \"47\u00f79=5, 2\"\n  }\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n\n  // Scope the search to this single cell's body so that identical text\n  // appearing elsewhere in the table (old values re-appear as new values\n  // for other cells) can never cause a cross-cell match.\n  const found = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\" in cell (${row},${col}), found ${found.items.length}.`\n    );\n  }\n\n  // insertText(..., replace) on the found range swaps only the text run's\n  // content, leaving the run's rPr (font/size) and the paragraph's pPr\n  // (alignment) untouched.\n  found.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Replaces the 25 arithmetic-answer text cells inside the single table in\n# the document with their new values, cell-by-cell, while preserving each\n# cell's existing run/paragraph formatting (font, size, alignment).\n#\n# Each entry records the 1-based (Row, Col) of a content-bearing table\n# cell (matching Table.Cell(row, col)), the text currently in that cell,\n# and the text it must become, taken straight from the authoritative diff.\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = '90\u00f73=30, 0'; NewText = '74\u00f78=9, 2' },\n    @{ Row = 1; Col = 2; OldText = '75\u00f72=37, 1'; NewText = '79\u00f73=26, 1' },\n    @{ Row = 1; Col = 3; OldText = '15\u00f79=1, 6'; NewText = '31\u00f76=5, 1' },\n    @{ Row = 1; Col = 4; OldText = '76\u00f74=19, 0'; NewText = '40\u00f75=8, 0' },\n    @{ Row = 1; Col = 5; OldText = '10\u00f73=3, 1'; NewText = '98\u00f74=24, 2' },\n    @{ Row = 5; Col = 1; OldText = '74\u00f78=9, 2'; NewText = '48\u00f74=12, 0' },\n    @{ Row = 5; Col = 2; OldText = '86\u00f73=28, 2'; NewText = '40\u00f74=10, 0' },\n    @{ Row = 5; Col = 3; OldText = '35\u00f75=7, 0'; NewText = '66\u00f74=16, 2' },\n    @{ Row = 5; Col = 4; OldText = '10\u00f76=1, 4'; NewText = '59\u00f77=8, 3' },\n    @{ Row = 5; Col = 5; OldText = '23\u00f75=4, 3'; NewText = '95\u00f75=19, 0' },\n    @{ Row = 9; Col = 1; OldText = '27\u00f72=13, 1'; NewText = '91\u00f73=30, 1' },\n    @{ Row = 9; Col = 2; OldText = '86\u00f78=10, 6'; NewText = '98\u00f72=49, 0' },\n    @{ Row = 9; Col = 3; OldText = '45\u00f76=7, 3'; NewText = '15\u00f78=1, 7' },\n    @{ Row = 9; Col = 4; OldText = '66\u00f74=16, 2'; NewText = '36\u00f79=4, 0' },\n    @{ Row = 9; Col = 5; OldText = '57\u00f75=11, 2'; NewText = '56\u00f74=14, 0' },\n    @{ Row = 13; Col = 1; OldText = '17\u00f79=1, 8'; NewText = '31\u00f75=6, 1' },\n    @{ Row = 13; Col = 2; OldText = '11\u00f74=2, 3'; NewText = '88\u00f76=14, 4' },\n    @{ Row = 13; Col = 3; OldText = '90\u00f72=45, 0'; NewText = '15\u00f73=5, 0' },\n    @{ Row = 13; Col = 4; OldText = '73\u00f73=24, 1'; NewText = '31\u00f72=15, 1' },\n    @{ Row = 13; Col = 5; OldText = '48\u00f73=16, 0'; NewText = '29\u00f74=7, 1' },\n    @{ Row = 17; Col = 1; OldText = '81\u00f79=9, 0'; NewText = '76\u00f72=38, 0' },\n    @{ Row = 17; Col = 2; OldText = '55\u00f77=7, 6'; NewText = '61\u00f75=12, 1' },\n    @{ Row = 17; Col = 3; OldText = '90\u00f75=18, 0'; NewText = '62\u00f75=12, 2' },\n    @{ Row = 17; Col = 4; OldText = '84\u00f75=16, 4'; NewText = '98\u00f72=49, 0' },\n    @{ Row = 17; Col = 5; OldText = '80\u00f78=10, 0'; NewText = '47\u00f79=5, 2' }\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nforeach ($r in $replacements) {\n    $cell = $tbl.Cell($r.Row, $r.Col)\n    $rng = $cell.Range\n\n    # Cell.Range.Text includes the trailing end-of-cell marker(s) (CR +\n    # BEL); strip those before comparing against the plain-text value we\n    # expect to find.\n    $current = $rng.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $r.OldText) {\n        throw \"Cell ($($r.Row),$($r.Col)) expected '$($r.OldText)' but found '$current'.\"\n    }\n\n    # Assigning Range.Text replaces only this range's content in place\n    # (it is confined to the cell, unlike a document-wide Find/Replace),\n    # and keeps the surrounding run formatting (rFonts/sz) and paragraph\n    # formatting (jc) untouched.\n    $rng.Text = $r.NewText\n}\n"}
